$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-order the ingredient lists ("Materias primas") for the existing recipes.
$ws.Range("C2").Value = "5.0-leche,1.0-vainilla,2.0-harinita,1.0-huevos,"
$ws.Range("C3").Value = "2.0-manzana,1.0-harinita,5.0-huevos,"
$ws.Range("C4").Value = "1.0-vainilla,5.0-harinita,2.0-huevos,"
$ws.Range("C5").Value = "5.0-harinita,5.0-huevos,"
$ws.Range("C6").Value = "1.0-crema,2.0-limon,4.0-harinita,5.0-huevos,5.0-merengue,"

# Add the new "asd" recipe row validating elaboration time.
$ws.Range("A9").Value = "asd"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "1.0-harinita,"
$ws.Range("D9").Value = 20
